$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = "Julien_ABrilliantPhoto_2022_05.jpg"
$ws.Range("B22").Value = "Samuel_ABrilliantPhoto_2022_07.jpg"
$ws.Range("B21").Value = "Rita_ABrilliantPhoto_2022_14.jpg"
$ws.Range("B19").Value = "Nora_ABrilliantPhoto_Feb2024_05.jpg"
$ws.Range("B18").Value = "Noemi_November2022_12_.jpg"
$ws.Range("B16").Value = "Lukas_ABrilliantPhoto_2022_04.jpg"
$ws.Range("B14").Value = "Jose_ABrilliantPhoto_2024_04.jpg"
$ws.Range("B13").Value = "Johannes_ABrilliantPhoto_2022_03.jpg"
$ws.Range("B12").Value = "Ivan_ABrilliantPhoto_2022_10.jpg"
$ws.Range("B11").Value = "Irmantas_ABrilliantPhoto_2024_09.jpg"
$ws.Range("B10").Value = "Flavie_ABrilliantPhoto_2022_06.jpg"
$ws.Range("B9").Value = "Dominique_ABrilliantPhoto_2022_06.jpg"
$ws.Range("B8").Value = "DanielaS_ABrilliantPhoto_2022_13.jpg"
$ws.Range("B7").Value = "DanielaM_ABrilliantPhoto_2022_06.jpg"
$ws.Range("B6").Value = "Christina_ABrilliantPhoto_2022_09.jpg"
$ws.Range("B5").Value = "Christian_ABrilliantPhoto_2022_19.jpg"
$ws.Range("B4").Value = "Julien_ABrilliantPhoto_2023_17.jpg"
$ws.Range("B3").Value = "Balduin_ABrilliantPhoto_2022_15.jpg"

$ws.Range("B15").Select()

$excel.ActiveWindow.Zoom = 221
